$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Move existing rows 2-21 down to rows 12-31 (write from bottom up to avoid overwrite issues)
$ws.Cells.Item(31, 1).Value = -4.002721786499023
$ws.Cells.Item(31, 2).Value = 1.022015571594239
$ws.Cells.Item(31, 3).Value = -0.0432633161544799
$ws.Cells.Item(30, 1).Value = 3.173869132995605
$ws.Cells.Item(30, 2).Value = -1.535699486732483
$ws.Cells.Item(30, 3).Value = -6.114311695098877
$ws.Cells.Item(29, 1).Value = -8.473310470581055
$ws.Cells.Item(29, 2).Value = -0.7327957153320312
$ws.Cells.Item(29, 3).Value = 5.200639724731445
$ws.Cells.Item(28, 1).Value = -2.372189998626709
$ws.Cells.Item(28, 2).Value = 1.225671410560608
$ws.Cells.Item(28, 3).Value = 2.504203796386719
$ws.Cells.Item(27, 1).Value = -1.961796522140503
$ws.Cells.Item(27, 2).Value = 1.68219518661499
$ws.Cells.Item(27, 3).Value = 1.394426345825195
$ws.Cells.Item(26, 1).Value = -3.03963303565979
$ws.Cells.Item(26, 2).Value = 1.802032470703125
$ws.Cells.Item(26, 3).Value = -2.227274417877197
$ws.Cells.Item(25, 1).Value = 0.8469958305358887
$ws.Cells.Item(25, 2).Value = -1.08077871799469
$ws.Cells.Item(25, 3).Value = 7.442714691162109
$ws.Cells.Item(24, 1).Value = -2.832679748535156
$ws.Cells.Item(24, 2).Value = 5.107204437255859
$ws.Cells.Item(24, 3).Value = -6.522222995758057
$ws.Cells.Item(23, 1).Value = -0.552617073059082
$ws.Cells.Item(23, 2).Value = 1.007189750671387
$ws.Cells.Item(23, 3).Value = -2.683732509613037
$ws.Cells.Item(22, 1).Value = -2.810617446899414
$ws.Cells.Item(22, 2).Value = 0.8466755151748657
$ws.Cells.Item(22, 3).Value = -0.6261429786682129
$ws.Cells.Item(21, 1).Value = -1.000519752502441
$ws.Cells.Item(21, 2).Value = -0.010628342628479
$ws.Cells.Item(21, 3).Value = -1.670511245727539
$ws.Cells.Item(20, 1).Value = 4.286171913146973
$ws.Cells.Item(20, 2).Value = 0.2758489847183227
$ws.Cells.Item(20, 3).Value = -4.509784698486328
$ws.Cells.Item(19, 1).Value = -11.09067344665527
$ws.Cells.Item(19, 2).Value = 1.405970811843872
$ws.Cells.Item(19, 3).Value = 10.02403450012207
$ws.Cells.Item(18, 1).Value = -2.353589773178101
$ws.Cells.Item(18, 2).Value = 0.5766786336898804
$ws.Cells.Item(18, 3).Value = 2.404436111450196
$ws.Cells.Item(17, 1).Value = -1.905292510986328
$ws.Cells.Item(17, 2).Value = 1.267569422721863
$ws.Cells.Item(17, 3).Value = 0.3008813858032226
$ws.Cells.Item(16, 1).Value = 0.1245284080505371
$ws.Cells.Item(16, 2).Value = 0.4134435057640075
$ws.Cells.Item(16, 3).Value = 2.055456638336182
$ws.Cells.Item(15, 1).Value = -0.6316938400268555
$ws.Cells.Item(15, 2).Value = 0.0533061251044273
$ws.Cells.Item(15, 3).Value = -1.823783159255981
$ws.Cells.Item(14, 1).Value = 0.7375173568725586
$ws.Cells.Item(14, 2).Value = -0.8549392819404602
$ws.Cells.Item(14, 3).Value = -2.997310400009156
$ws.Cells.Item(13, 1).Value = -0.1681756973266601
$ws.Cells.Item(13, 2).Value = -0.045459896326065
$ws.Cells.Item(13, 3).Value = 0.3079473972320556
$ws.Cells.Item(12, 1).Value = -0.1584005355834961
$ws.Cells.Item(12, 2).Value = 0.0559865832328796
$ws.Cells.Item(12, 3).Value = -0.2031860947608947

# Step 2: Write new rows 2-11
$ws.Cells.Item(2, 1).Value = -0.1118526458740234
$ws.Cells.Item(2, 2).Value = 0.0269185900688171
$ws.Cells.Item(2, 3).Value = 0.0618541836738586
$ws.Cells.Item(3, 1).Value = -0.188694953918457
$ws.Cells.Item(3, 2).Value = -0.0127399563789367
$ws.Cells.Item(3, 3).Value = 0.0153613984584808
$ws.Cells.Item(4, 1).Value = -0.0261173248291015
$ws.Cells.Item(4, 2).Value = -0.1474769711494445
$ws.Cells.Item(4, 3).Value = 0.0655251443386077
$ws.Cells.Item(5, 1).Value = -0.1960973739624023
$ws.Cells.Item(5, 2).Value = 0.0549294650554657
$ws.Cells.Item(5, 3).Value = 0.0360765755176544
$ws.Cells.Item(6, 1).Value = -0.0661020278930664
$ws.Cells.Item(6, 2).Value = -0.1787786185741424
$ws.Cells.Item(6, 3).Value = 0.0745508223772049
$ws.Cells.Item(7, 1).Value = 0.0234136581420898
$ws.Cells.Item(7, 2).Value = 0.0270741879940032
$ws.Cells.Item(7, 3).Value = 0.2239813506603241
$ws.Cells.Item(8, 1).Value = 0.11651611328125
$ws.Cells.Item(8, 2).Value = -0.4856438636779785
$ws.Cells.Item(8, 3).Value = 0.5658785104751587
$ws.Cells.Item(9, 1).Value = 0.0557413101196289
$ws.Cells.Item(9, 2).Value = 0.3574482798576355
$ws.Cells.Item(9, 3).Value = 0.2321825623512268
$ws.Cells.Item(10, 1).Value = 0.3619680404663086
$ws.Cells.Item(10, 2).Value = 0.0124948024749755
$ws.Cells.Item(10, 3).Value = 0.3587799966335296
$ws.Cells.Item(11, 1).Value = -0.2529764175415039
$ws.Cells.Item(11, 2).Value = 0.1160029470920562
$ws.Cells.Item(11, 3).Value = -0.0988222360610961

Write-Host "Done"
